$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "yearly period" headers in row 8 and row 24 ---
# The oldest period (1396/12) is dropped and a new one (1401/12) is
# appended, so every existing period shifts one column to the left
# (E<-F, F<-G, G<-H, H<-I) and the newest period (1401/12) lands in I.
$periods = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)

for ($i = 0; $i -lt 5; $i++) {
    $col = 5 + $i   # E=5 .. I=9
    $ws.Cells.Item(8, $col).Value = $periods[$i]
    $ws.Cells.Item(24, $col).Value = $periods[$i]
}

# --- Update the numeric data rows: shift existing values left one
# column and fill the new rightmost column (I) with the newly
# reported figure. ---
$dataRows = @{
    13 = @(12843, 26668, 59710, 83001, 3970)
    14 = @(1310, 1147, 2303, 6734, 5734)
    15 = @(305, 468, 129, 882, 544)
    16 = @(2968, 3965, 6814, 10351, 6066)
    17 = @(48783, 56867, 76397, 108762, 202753)
    19 = @(96746, 139316, 188819, 246506, 377340)
    20 = @(162955, 228431, 334172, 456236, 596407)
    26 = @(82, 60, 58, 58, 55)
    27 = @(359, 388, 455, 455, 533)
}

foreach ($row in $dataRows.Keys) {
    $values = $dataRows[$row]
    for ($i = 0; $i -lt 5; $i++) {
        $col = 5 + $i   # E=5 .. I=9
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
